# Auto-generated edit script applying numeric corrections to the Ifrit_Profits workbook
# (scheduled runner update touching currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.83333
$ws.Range("I9").Value = 141.25
$ws.Range("J9").Value = 152
$ws.Range("K9").Value = 141.25
$ws.Range("L9").Value = 152
$ws.Range("M9").Value = 27.75
$ws.Range("N9").Value = -490
$ws.Range("H111").Value = 1475.1538
$ws.Range("I111").Value = 1546.8422
$ws.Range("J111").Value = 1280.5714
$ws.Range("K111").Value = 4640.5266
$ws.Range("L111").Value = 3841.7142
$ws.Range("M111").Value = -1573.5266
$ws.Range("N111").Value = -9975.7142
$ws.Range("H137").Value = 27779524
$ws.Range("I137").Value = 1380.25
$ws.Range("J137").Value = 62502204
$ws.Range("K137").Value = 4140.75
$ws.Range("L137").Value = 187506612
$ws.Range("M137").Value = -1590.75
$ws.Range("N137").Value = -187511712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5534.622
$ws.Range("I32").Value = 4563.684
$ws.Range("K32").Value = 4563.684
$ws.Range("M32").Value = -4276.684
$ws.Range("H74").Value = 8346371.5
$ws.Range("I74").Value = 50001356
$ws.Range("J74").Value = 15374.5
$ws.Range("K74").Value = 50001356
$ws.Range("L74").Value = 15374.5
$ws.Range("M74").Value = -50000482
$ws.Range("N74").Value = -17122.5
$ws.Range("H77").Value = 8346371.5
$ws.Range("I77").Value = 50001356
$ws.Range("J77").Value = 15374.5
$ws.Range("K77").Value = 250006780
$ws.Range("L77").Value = 76872.5
$ws.Range("M77").Value = -250002412
$ws.Range("N77").Value = -85608.5
$ws.Range("H110").Value = 1184.2941
$ws.Range("I110").Value = 911.4167
$ws.Range("J110").Value = 1839.2
$ws.Range("K110").Value = 911.4167
$ws.Range("L110").Value = 1839.2
$ws.Range("M110").Value = 1133.5833
$ws.Range("N110").Value = -5929.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1507.6923
$ws.Range("I105").Value = 1460
$ws.Range("J105").Value = 1666.6666
$ws.Range("K105").Value = 1460
$ws.Range("L105").Value = 1666.6666
$ws.Range("M105").Value = 287
$ws.Range("N105").Value = -5160.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 35714496
$ws.Range("I33").Value = 38461750
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 230770500
$ws.Range("L33").Value = 720
$ws.Range("M33").Value = -230770217
$ws.Range("N33").Value = -1286
$ws.Range("H44").Value = 881.25
$ws.Range("I44").Value = 387.5
$ws.Range("J44").Value = 1375
$ws.Range("K44").Value = 1162.5
$ws.Range("L44").Value = 4125
$ws.Range("M44").Value = -764.5
$ws.Range("N44").Value = -4921
$ws.Range("H47").Value = 273.25
$ws.Range("I47").Value = 197.66667
$ws.Range("J47").Value = 500
$ws.Range("K47").Value = 593.00001
$ws.Range("L47").Value = 1500
$ws.Range("M47").Value = -162.00001
$ws.Range("N47").Value = -2362
$ws.Range("H63").Value = 2006
$ws.Range("I63").Value = 1012
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3036
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -2287
$ws.Range("N63").Value = -10498
$ws.Range("H64").Value = 1401470.2
$ws.Range("I64").Value = 955.75
$ws.Range("J64").Value = 2335146.8
$ws.Range("K64").Value = 2867.25
$ws.Range("L64").Value = 7005440.399999999
$ws.Range("M64").Value = -2597.25
$ws.Range("N64").Value = -7005980.399999999
$ws.Range("H66").Value = 2006
$ws.Range("I66").Value = 1012
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 9108
$ws.Range("L66").Value = 27000
$ws.Range("M66").Value = -5364
$ws.Range("N66").Value = -34488
$ws.Range("H67").Value = 1401470.2
$ws.Range("I67").Value = 955.75
$ws.Range("J67").Value = 2335146.8
$ws.Range("K67").Value = 2867.25
$ws.Range("L67").Value = 7005440.399999999
$ws.Range("M67").Value = -1931.25
$ws.Range("N67").Value = -7007312.399999999
$ws.Range("H81").Value = 1135.8572
$ws.Range("J81").Value = 1200.1666
$ws.Range("L81").Value = 3600.4998
$ws.Range("N81").Value = -5846.4998
$ws.Range("H84").Value = 1135.8572
$ws.Range("J84").Value = 1200.1666
$ws.Range("L84").Value = 10801.4994
$ws.Range("N84").Value = -22033.4994
$ws.Range("H92").Value = 503
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 503
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1509
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4005
$ws.Range("H94").Value = 1390
$ws.Range("I94").Value = 1390
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4170
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3494
$ws.Range("N94").ClearContents()
$ws.Range("H95").Value = 3000
$ws.Range("J95").Value = 3000
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -13118
$ws.Range("H97").Value = 536.2105
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H98").Value = 273
$ws.Range("I98").Value = 263.25
$ws.Range("J98").Value = 279.5
$ws.Range("K98").Value = 789.75
$ws.Range("L98").Value = 838.5
$ws.Range("M98").Value = 708.25
$ws.Range("N98").Value = -3834.5
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 4500
$ws.Range("M102").Value = -2066
$ws.Range("H103").Value = 2428905.5
$ws.Range("I103").Value = 3091113.5
$ws.Range("J103").Value = 809.3333
$ws.Range("K103").Value = 9273340.5
$ws.Range("L103").Value = 2427.9999
$ws.Range("M103").Value = -9272461.5
$ws.Range("N103").Value = -4185.9999
$ws.Range("H131").Value = 1972.1111
$ws.Range("I131").Value = 7682.5713
$ws.Range("J131").Value = 1490.506
$ws.Range("K131").Value = 23047.7139
$ws.Range("L131").Value = 4471.518
$ws.Range("M131").Value = -18007.7139
$ws.Range("N131").Value = -14551.518
$ws.Range("H137").Value = 5718.375
$ws.Range("I137").Value = 3778
$ws.Range("J137").Value = 6077.7036
$ws.Range("K137").Value = 11334
$ws.Range("L137").Value = 18233.1108
$ws.Range("M137").Value = -6234
$ws.Range("N137").Value = -28433.1108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 16918.45
$ws.Range("I122").Value = 2818.1667
$ws.Range("J122").Value = 39991.637
$ws.Range("K122").Value = 8454.500100000001
$ws.Range("L122").Value = 119974.911
$ws.Range("M122").Value = -6004.500100000001
$ws.Range("N122").Value = -124874.911
$ws.Range("H126").Value = 2449.9
$ws.Range("J126").Value = 2999.8
$ws.Range("L126").Value = 8999.400000000001
$ws.Range("N126").Value = -13939.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20167.166
$ws.Range("I132").Value = 23000.8
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 69002.39999999999
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -66472.39999999999
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716
$ws.Range("H132").Value = 5484.448
$ws.Range("I132").Value = 6965.579
$ws.Range("J132").Value = 2670.3
$ws.Range("K132").Value = 20896.737
$ws.Range("L132").Value = 8010.900000000001
$ws.Range("M132").Value = -18366.737
$ws.Range("N132").Value = -13070.9
